$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F108
$ws.Range("F2").Value = "2021-10-05 13:40:17.492157"
$ws.Range("F3").Value = "2021-10-05 13:40:17.492170"
$ws.Range("F4").Value = "2021-10-05 13:40:17.492175"
$ws.Range("F5").Value = "2021-10-05 13:40:17.492178"
$ws.Range("F6").Value = "2021-10-05 13:40:17.492181"
$ws.Range("F7").Value = "2021-10-05 13:40:17.492184"
$ws.Range("F8").Value = "2021-10-05 13:40:17.492187"
$ws.Range("F9").Value = "2021-10-05 13:40:17.492190"
$ws.Range("F10").Value = "2021-10-05 13:40:17.492194"
$ws.Range("F11").Value = "2021-10-05 13:40:17.492197"
$ws.Range("F12").Value = "2021-10-05 13:40:17.492200"
$ws.Range("F13").Value = "2021-10-05 13:40:17.492203"
$ws.Range("F14").Value = "2021-10-05 13:40:17.492206"
$ws.Range("F15").Value = "2021-10-05 13:40:17.492209"
$ws.Range("F16").Value = "2021-10-05 13:40:17.492212"
$ws.Range("F17").Value = "2021-10-05 13:40:17.492214"
$ws.Range("F18").Value = "2021-10-05 13:40:17.492218"
$ws.Range("F19").Value = "2021-10-05 13:40:17.492221"
$ws.Range("F20").Value = "2021-10-05 13:40:17.492224"
$ws.Range("F21").Value = "2021-10-05 13:40:17.492227"
$ws.Range("F22").Value = "2021-10-05 13:40:17.492230"
$ws.Range("F23").Value = "2021-10-05 13:40:17.492233"
$ws.Range("F24").Value = "2021-10-05 13:40:17.492236"
$ws.Range("F25").Value = "2021-10-05 13:40:17.492239"
$ws.Range("F26").Value = "2021-10-05 13:40:17.492242"
$ws.Range("F27").Value = "2021-10-05 13:40:17.492245"
$ws.Range("F28").Value = "2021-10-05 13:40:17.492248"
$ws.Range("F29").Value = "2021-10-05 13:40:17.492251"
$ws.Range("F30").Value = "2021-10-05 13:40:17.492254"
$ws.Range("F31").Value = "2021-10-05 13:40:17.492257"
$ws.Range("F32").Value = "2021-10-05 13:40:17.492260"
$ws.Range("F33").Value = "2021-10-05 13:40:17.492263"
$ws.Range("F34").Value = "2021-10-05 13:40:17.492267"
$ws.Range("F35").Value = "2021-10-05 13:40:17.492270"
$ws.Range("F36").Value = "2021-10-05 13:40:17.492273"
$ws.Range("F37").Value = "2021-10-05 13:40:17.492276"
$ws.Range("F38").Value = "2021-10-05 13:40:17.492279"
$ws.Range("F39").Value = "2021-10-05 13:40:17.492282"
$ws.Range("F40").Value = "2021-10-05 13:40:17.492284"
$ws.Range("F41").Value = "2021-10-05 13:40:17.492287"
$ws.Range("F42").Value = "2021-10-05 13:40:17.492291"
$ws.Range("F43").Value = "2021-10-05 13:40:17.492294"
$ws.Range("F44").Value = "2021-10-05 13:40:17.492297"
$ws.Range("F45").Value = "2021-10-05 13:40:17.492300"
$ws.Range("F46").Value = "2021-10-05 13:40:17.492303"
$ws.Range("F47").Value = "2021-10-05 13:40:17.492306"
$ws.Range("F48").Value = "2021-10-05 13:40:17.492309"
$ws.Range("F49").Value = "2021-10-05 13:40:17.492312"
$ws.Range("F50").Value = "2021-10-05 13:40:17.492315"
$ws.Range("F51").Value = "2021-10-05 13:40:17.492318"
$ws.Range("F52").Value = "2021-10-05 13:40:17.492320"
$ws.Range("F53").Value = "2021-10-05 13:40:17.492323"
$ws.Range("F54").Value = "2021-10-05 13:40:17.492327"
$ws.Range("F55").Value = "2021-10-05 13:40:17.492330"
$ws.Range("F56").Value = "2021-10-05 13:40:17.492333"
$ws.Range("F57").Value = "2021-10-05 13:40:17.492336"
$ws.Range("F58").Value = "2021-10-05 13:40:17.492339"
$ws.Range("F59").Value = "2021-10-05 13:40:17.492342"
$ws.Range("F60").Value = "2021-10-05 13:40:17.492345"
$ws.Range("F61").Value = "2021-10-05 13:40:17.492347"
$ws.Range("F62").Value = "2021-10-05 13:40:17.492350"
$ws.Range("F63").Value = "2021-10-05 13:40:17.492353"
$ws.Range("F64").Value = "2021-10-05 13:40:17.492360"
$ws.Range("F65").Value = "2021-10-05 13:40:17.492364"
$ws.Range("F66").Value = "2021-10-05 13:40:17.492368"
$ws.Range("F67").Value = "2021-10-05 13:40:17.492371"
$ws.Range("F68").Value = "2021-10-05 13:40:17.492389"
$ws.Range("F69").Value = "2021-10-05 13:40:17.492391"
$ws.Range("F70").Value = "2021-10-05 13:40:17.492394"
$ws.Range("F71").Value = "2021-10-05 13:40:17.492396"
$ws.Range("F72").Value = "2021-10-05 13:40:17.492399"
$ws.Range("F73").Value = "2021-10-05 13:40:17.492402"
$ws.Range("F74").Value = "2021-10-05 13:40:17.492404"
$ws.Range("F75").Value = "2021-10-05 13:40:17.492407"
$ws.Range("F76").Value = "2021-10-05 13:40:17.492409"
$ws.Range("F77").Value = "2021-10-05 13:40:17.492412"
$ws.Range("F78").Value = "2021-10-05 13:40:17.492416"
$ws.Range("F79").Value = "2021-10-05 13:40:17.492420"
$ws.Range("F80").Value = "2021-10-05 13:40:17.492422"
$ws.Range("F81").Value = "2021-10-05 13:40:17.492425"
$ws.Range("F82").Value = "2021-10-05 13:40:17.492427"
$ws.Range("F83").Value = "2021-10-05 13:40:17.492430"
$ws.Range("F84").Value = "2021-10-05 13:40:17.492433"
$ws.Range("F85").Value = "2021-10-05 13:40:17.492436"
$ws.Range("F86").Value = "2021-10-05 13:40:17.492438"
$ws.Range("F87").Value = "2021-10-05 13:40:17.492441"
$ws.Range("F88").Value = "2021-10-05 13:40:17.492443"
$ws.Range("F89").Value = "2021-10-05 13:40:17.492446"
$ws.Range("F90").Value = "2021-10-05 13:40:17.492448"
$ws.Range("F91").Value = "2021-10-05 13:40:17.492451"
$ws.Range("F92").Value = "2021-10-05 13:40:17.492453"
$ws.Range("F93").Value = "2021-10-05 13:40:17.492456"
$ws.Range("F94").Value = "2021-10-05 13:40:17.492460"
$ws.Range("F95").Value = "2021-10-05 13:40:17.492462"
$ws.Range("F96").Value = "2021-10-05 13:40:17.492465"
$ws.Range("F97").Value = "2021-10-05 13:40:17.492468"
$ws.Range("F98").Value = "2021-10-05 13:40:17.492470"
$ws.Range("F99").Value = "2021-10-05 13:40:17.492473"
$ws.Range("F100").Value = "2021-10-05 13:40:17.492475"
$ws.Range("F101").Value = "2021-10-05 13:40:17.492478"
$ws.Range("F102").Value = "2021-10-05 13:40:17.492481"
$ws.Range("F103").Value = "2021-10-05 13:40:17.492483"
$ws.Range("F104").Value = "2021-10-05 13:40:17.492486"
$ws.Range("F105").Value = "2021-10-05 13:40:17.492488"
$ws.Range("F106").Value = "2021-10-05 13:40:17.492491"
$ws.Range("F107").Value = "2021-10-05 13:40:17.492494"
$ws.Range("F108").Value = "2021-10-05 13:40:17.492496"

# Copy formatting from E1 (existing header style) onto F1 so it matches the other headers
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

